$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text formatting instead of
# being auto-converted to a number by Excel when the new values are plain
# decimal numbers (e.g. "584.41").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.491.72"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "3.464.17"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "584.41"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").Value = "177.79"
$ws.Range("E6").Value = "  +1.81%  "
$ws.Range("D7").Value = "0.629"
$ws.Range("E7").Value = "  +5.91%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "3.463.06"
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").Value = "0.134"
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("D11").Value = "6.97"
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("D13").Value = "4.066.34"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("D15").Value = "30.20"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").Value = "66.379.65"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "3.465.88"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").Value = "5.97"
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("D20").Value = "13.86"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").Value = "372.49"
$ws.Range("E21").Value = "  -2.28%  "
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("D23").Value = "73.39"
$ws.Range("E23").Value = "  +1.61%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -1.44%  "
$ws.Range("E26").Value = "  +5.26%  "
$ws.Range("D27").Value = "10.06"
$ws.Range("E27").Value = "  +2.35%  "
$ws.Range("E28").Value = "  +3.48%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "5.97"
$ws.Range("E30").Value = "  +1.52%  "
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").Value = "23.75"
$ws.Range("E32").Value = "  -3.23%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "7.07"
$ws.Range("E34").Value = "  -2.45%  "
$ws.Range("E35").Value = "  -4.74%  "
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("D37").Value = "160.93"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "0.887"
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("E39").Value = "  -5.11%  "
$ws.Range("E40").Value = "  +1.91%  "
$ws.Range("D41").Value = "2.812.65"
$ws.Range("E41").Value = "  +3.32%  "
$ws.Range("D42").Value = "4.53"
$ws.Range("E42").Value = "  +0.90%  "
$ws.Range("D43").Value = "2.59"
$ws.Range("E43").Value = "  +2.77%  "
$ws.Range("D44").Value = "6.49"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("D46").Value = "25.29"
$ws.Range("E46").Value = "  +1.30%  "
$ws.Range("D47").Value = "342.47"
$ws.Range("E47").Value = "  +5.90%  "
$ws.Range("D48").Value = "40.07"
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("E50").Value = "  +3.02%  "
$ws.Range("E51").Value = "  -0.87%  "
